$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data (row 17), matching the date formatting used by the rows above it
$ws.Range("A16").Copy()
$ws.Range("A17").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("A17").Value = (Get-Date -Year 2025 -Month 3 -Day 21 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B17").Value = "1 hour"
$ws.Range("C17").Value = "Update data and plots"

# Update the active selection to reflect where the user ended up after editing
$ws.Range("C18").Select()
